# "Minor tweaks - getting close to final!"
#
# Slides 23 ("SQL Injection" title slide) and 24 ("SQL Injection" content
# slide) are being hidden from the live run-through, and both get a slow
# 2-second transition applied.

$p = $ppt.ActivePresentation

$slideIndexes = @(23, 24)

foreach ($idx in $slideIndexes) {
    $s = $p.Slides.Item($idx)

    # Set the transition duration/speed before hiding so the written
    # <p:transition> element carries both the legacy speed attribute and
    # the precise (p14) duration in milliseconds.
    $s.SlideShowTransition.Duration = 2
    $s.SlideShowTransition.Speed = 1   # ppTransitionSpeedSlow
    $s.SlideShowTransition.Hidden = $true
}
